$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Update the handlebars note text to expose the data model (data.refrigerator_id)
$ws.Range("D2").Value = "Refrigerator id: {{data.refrigerator_id}}"

# Make the survey sheet the active tab/selection
$ws.Activate()
$ws.Range("D2").Select()
